$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 67
$ws.Range("A67").Value = 111870990
$ws.Range("B67").Value = 90666
$ws.Range("D67").Value = "LC"
$ws.Range("E67").Value = 4364
$ws.Range("F67").Value = "Dropptaggsvamp"
$ws.Range("G67").Value = "Hydnellum ferrugineum"
$ws.Range("H67").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q67").Value = 590569.8478412227
$ws.Range("R67").Value = 7040376.109235858

# Row 68
$ws.Range("A68").Value = 111881322
$ws.Range("B68").Value = 56414
$ws.Range("D68").Value = "NT"
$ws.Range("E68").Value = 100049
$ws.Range("F68").Value = "Spillkråka"
$ws.Range("G68").Value = "Dryocopus martius"
$ws.Range("H68").Value = "(Linnaeus, 1758)"
$ws.Range("P68").Value = "Valforsen, Ång"
$ws.Range("Q68").Value = 590615.1562677342
$ws.Range("R68").Value = 7040278.573758457
$ws.Range("K68").ClearContents()
$ws.Range("M68").Value = "gammalt bo"

# Row 69
$ws.Range("A69").Value = 111870139
$ws.Range("B69").Value = 89845
$ws.Range("D69").Value = "VU"
$ws.Range("E69").Value = 1209
$ws.Range("F69").Value = "Rynkskinn"
$ws.Range("G69").Value = "Phlebia centrifuga"
$ws.Range("H69").Value = "P.Karst."
$ws.Range("P69").Value = "Valforsen (Valforsen), Ång"
$ws.Range("Q69").Value = 590710.4131779457
$ws.Range("R69").Value = 7040581.765558361
$ws.Range("K69").Value = ""

# Row 70
$ws.Range("A70").Value = 111881310
$ws.Range("B70").Value = 89425
$ws.Range("E70").Value = 5442
$ws.Range("F70").Value = "Tallticka"
$ws.Range("G70").Value = "Porodaedalea pini"
$ws.Range("H70").Value = "(Brot.) Murrill"
$ws.Range("P70").Value = "Valforsen, Ång"
$ws.Range("Q70").Value = 590738.9206925276
$ws.Range("R70").Value = 7040524.002523924
$ws.Range("K70").ClearContents()

# Row 71
$ws.Range("A71").Value = 111871585
$ws.Range("B71").Value = 89405
$ws.Range("D71").Value = "NT"
$ws.Range("E71").Value = 1202
$ws.Range("F71").Value = "Ullticka"
$ws.Range("G71").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H71").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q71").Value = 590630.2636057099
$ws.Range("R71").Value = 7040266.929520278

# Row 72
$ws.Range("A72").Value = 111870830
$ws.Range("B72").Value = 90678
$ws.Range("D72").Value = "LC"
$ws.Range("E72").Value = 4366
$ws.Range("F72").Value = "Skarp dropptaggsvamp"
$ws.Range("G72").Value = "Hydnellum peckii"
$ws.Range("H72").Value = "Banker"
$ws.Range("Q72").Value = 590558.4251677697
$ws.Range("R72").Value = 7040399.931061053

# Row 73
$ws.Range("A73").Value = 111870127
$ws.Range("B73").Value = 89405
$ws.Range("E73").Value = 1202
$ws.Range("F73").Value = "Ullticka"
$ws.Range("G73").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H73").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P73").Value = "Valforsen (Valforsen), Ång"
$ws.Range("Q73").Value = 590710.4131779457
$ws.Range("R73").Value = 7040581.765558361
$ws.Range("K73").Value = ""
$ws.Range("M73").ClearContents()
